# Weekly update of Pepino dulce (Femacal de La Calera) price rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44238
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 13000
$ws.Range("M2").Value = 13000
$ws.Range("P2").Value = 722

# Row 3
$ws.Range("D3").Value = 44238
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 11000
$ws.Range("P3").Value = 611

# Row 4
$ws.Range("D4").Value = 44424
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 18000
$ws.Range("N4").Value = "`$/caja 15 kilos"
$ws.Range("P4").Value = 1200
$ws.Range("Q4").Value = 15

# Row 5
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("P5").Value = 800

# Row 6
$ws.Range("D6").Value = 44235
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 14000
$ws.Range("N6").Value = "`$/bandeja 18 kilos"
$ws.Range("P6").Value = 778
$ws.Range("Q6").Value = 18

# Row 7
$ws.Range("D7").Value = 44235
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 12000
$ws.Range("P7").Value = 667

# Row 8
$ws.Range("D8").Value = 44235
$ws.Range("I8").Value = "Tercera"
$ws.Range("J8").Value = 60

# Row 9
$ws.Range("D9").Value = 44242
$ws.Range("J9").Value = 60

# Row 10
$ws.Range("D10").Value = 44242
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("P10").Value = 556
